# Add a new "OECD Air Emissions" column (column V) to the
# "Datasets and Years" sheet, matching the formatting already used by the
# neighbouring header/data cells (e.g. column B / column T), then update
# the view (scroll + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header cell (row 3) ------------------------------------------------
$ws.Range("V3").Value = "OECD Air Emissions"
$ws.Range("B3").Copy()
$ws.Range("V3").PasteSpecial($xlPasteFormats)

# --- Checkmarks for rows 4 through 31 -----------------------------------
$ws.Range("B4").Copy()
for ($r = 4; $r -le 31; $r++) {
    $cell = $ws.Range("V$r")
    $cell.Value = [char]0x2713
    $cell.PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

# --- View: scroll so row 8 is at the top, select A31 --------------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A31").Select()
